# Insert two new data rows (uva "Superior Seedless" Primera/Segunda entries
# dated 2022-03-18 [serial 44623], Región de O'Higgins) right after the
# header-adjacent data block, at sheet rows 26-27. This pushes every
# existing row from 26 onward down by two (old row 26 -> new row 28, etc.),
# matching the target diff (dimension grows from T122 to T124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two blank rows at 26 (Excel shifts 26.. down to 28..).
$ws.Rows("26:27").Insert()

# --- New row 26 --------------------------------------------------------
$ws.Cells.Item(26, 1).Value2  = 7
$ws.Cells.Item(26, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value2  = "Ñuble"
$ws.Cells.Item(26, 4).Value2  = 44623
$ws.Cells.Item(26, 5).Value2  = 16
$ws.Cells.Item(26, 6).Value2  = "Fruta"
$ws.Cells.Item(26, 7).Value2  = 100109
$ws.Cells.Item(26, 8).Value2  = "Uva"
$ws.Cells.Item(26, 9).Value2  = 100109001
$ws.Cells.Item(26, 10).Value2 = "Uva"
$ws.Cells.Item(26, 11).Value2 = "Superior Seedless"
$ws.Cells.Item(26, 12).Value2 = "Primera"
$ws.Cells.Item(26, 13).Value2 = 160
$ws.Cells.Item(26, 14).Value2 = 10000
$ws.Cells.Item(26, 15).Value2 = 11000
$ws.Cells.Item(26, 16).Value2 = 10500
$ws.Cells.Item(26, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(26, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(26, 19).Value2 = 583
$ws.Cells.Item(26, 20).Value2 = 18

# --- New row 27 --------------------------------------------------------
$ws.Cells.Item(27, 1).Value2  = 7
$ws.Cells.Item(27, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(27, 3).Value2  = "Ñuble"
$ws.Cells.Item(27, 4).Value2  = 44623
$ws.Cells.Item(27, 5).Value2  = 16
$ws.Cells.Item(27, 6).Value2  = "Fruta"
$ws.Cells.Item(27, 7).Value2  = 100109
$ws.Cells.Item(27, 8).Value2  = "Uva"
$ws.Cells.Item(27, 9).Value2  = 100109001
$ws.Cells.Item(27, 10).Value2 = "Uva"
$ws.Cells.Item(27, 11).Value2 = "Superior Seedless"
$ws.Cells.Item(27, 12).Value2 = "Segunda"
$ws.Cells.Item(27, 13).Value2 = 80
$ws.Cells.Item(27, 14).Value2 = 9000
$ws.Cells.Item(27, 15).Value2 = 9000
$ws.Cells.Item(27, 16).Value2 = 9000
$ws.Cells.Item(27, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(27, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(27, 19).Value2 = 500
$ws.Cells.Item(27, 20).Value2 = 18

# Apply the same date display format (style of column D) used by the rest
# of the "Fecha" column to the two newly-populated date cells.
$ws.Range("D26").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("D27").NumberFormat = $ws.Range("D28").NumberFormat
